$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00228310502283105
$ws.Range("C2").Value = 0.0045662100456621
$ws.Range("D2").Value = 0.00228310502283105
$ws.Range("E2").Value = 0.00228310502283105
$ws.Range("F2").Value = 0.045662100456621
$ws.Range("I2").Value = 0.0045662100456621
$ws.Range("J2").Value = 0.0045662100456621
$ws.Range("K2").Value = 0.0114155251141553
$ws.Range("L2").Value = 0.00684931506849315
$ws.Range("M2").Value = 0.0045662100456621
$ws.Range("O2").Value = 0.0182648401826484
$ws.Range("P2").Value = 0.0045662100456621
$ws.Range("Q2").Value = 0.00684931506849315
$ws.Range("R2").Value = 0.988584474885845
$ws.Range("T2").Value = 0.0045662100456621
$ws.Range("V2").Value = 0.0045662100456621
$ws.Range("W2").Value = 0.00684931506849315
$ws.Range("X2").Value = 0.0136986301369863
$ws.Range("B3").Value = 0.00684931506849315
$ws.Range("C3").Value = 0.972602739726027
$ws.Range("D3").Value = 0.977168949771689
$ws.Range("E3").Value = 0.984018264840183
$ws.Range("F3").Value = 0.0045662100456621
$ws.Range("G3").Value = 0.0159817351598174
$ws.Range("H3").Value = 0.974885844748858
$ws.Range("I3").Value = 0.00228310502283105
$ws.Range("J3").Value = 0.0045662100456621
$ws.Range("K3").Value = 0.0045662100456621
$ws.Range("M3").Value = 0.0091324200913242
$ws.Range("N3").Value = 0.00228310502283105
$ws.Range("R3").Value = 0.00684931506849315
$ws.Range("T3").Value = 0.0114155251141553
$ws.Range("U3").Value = 0.974885844748858
$ws.Range("W3").Value = 0.00684931506849315
$ws.Range("X3").Value = 0.00228310502283105
$ws.Range("B4").Value = 0.970319634703196
$ws.Range("D4").Value = 0.0159817351598174
$ws.Range("E4").Value = 0.0045662100456621
$ws.Range("F4").Value = 0.949771689497717
$ws.Range("G4").Value = 0.00228310502283105
$ws.Range("H4").Value = 0.0114155251141553
$ws.Range("I4").Value = 0.00684931506849315
$ws.Range("J4").Value = 0.988584474885845
$ws.Range("K4").Value = 0.0045662100456621
$ws.Range("L4").Value = 0.00684931506849315
$ws.Range("M4").Value = 0.974885844748858
$ws.Range("N4").Value = 0.0091324200913242
$ws.Range("O4").Value = 0.970319634703196
$ws.Range("P4").Value = 0.0091324200913242
$ws.Range("Q4").Value = 0.990867579908676
$ws.Range("U4").Value = 0.00684931506849315
$ws.Range("V4").Value = 0.0136986301369863
$ws.Range("W4").Value = 0.986301369863014
$ws.Range("X4").Value = 0.977168949771689
$ws.Range("B5").Value = 0.0114155251141553
$ws.Range("C5").Value = 0.0228310502283105
$ws.Range("D5").Value = 0.0045662100456621
$ws.Range("E5").Value = 0.0091324200913242
$ws.Range("G5").Value = 0.981735159817352
$ws.Range("H5").Value = 0.0136986301369863
$ws.Range("I5").Value = 0.986301369863014
$ws.Range("J5").Value = 0.00228310502283105
$ws.Range("K5").Value = 0.979452054794521
$ws.Range("L5").Value = 0.986301369863014
$ws.Range("M5").Value = 0.0114155251141553
$ws.Range("N5").Value = 0.988584474885845
$ws.Range("O5").Value = 0.0114155251141553
$ws.Range("P5").Value = 0.986301369863014
$ws.Range("Q5").Value = 0.00228310502283105
$ws.Range("R5").Value = 0.0045662100456621
$ws.Range("T5").Value = 0.984018264840183
$ws.Range("U5").Value = 0.00684931506849315
$ws.Range("V5").Value = 0.981735159817352
$ws.Range("X5").Value = 0.00684931506849315
